# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) values for the second data row
# (row 3) on both the zh-cn and de-de localization status sheets, as a
# new handback report run produced later timestamps for that file pair.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Sheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-21 03:23:39"
$wsZhCn.Range("H3").Value = "2016-03-21 03:24:20"

$wsDeDe = $wb.Sheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-21 03:23:47"
$wsDeDe.Range("H3").Value = "2016-03-21 03:24:38"
